$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update product name (B1) and short name (B2) on the input sheet
$wsInput.Range("B1").Value = "2430-RBI-EI-DB-DL-REC-NOCOM-RNI-CTPD-SAR-MD-TR-1-DATE-VAR-INST-1st"
$wsInput.Range("B2").Value = "243e"

# Mirror the new product name on the output sheet, and match its styling
$wsOutput.Range("B1").Value = "2430-RBI-EI-DB-DL-REC-NOCOM-RNI-CTPD-SAR-MD-TR-1-DATE-VAR-INST-1st"
$wsOutput.Range("B1").Font.Bold = $false

# Reset the input sheet's selection/scroll position back to the top
$wsInput.Activate()
$wsInput.Range("B1").Select()
$excel.ActiveWindow.ScrollRow = 1
